$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.943.07"
$ws.Range("E2").Value = "  -4.18%  "

# Row 3
$ws.Range("D3").Value = "2.230.69"
$ws.Range("E3").Value = "  -4.91%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "230.51"
$ws.Range("E5").Value = "  -3.79%  "

# Row 6
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  -6.70%  "

# Row 7
$ws.Range("D7").Value = "69.68"
$ws.Range("E7").Value = "  -4.50%  "

# Row 8
$ws.Range("E8").Value = "  +0.14%  "

# Row 9
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  -6.71%  "

# Row 10
$ws.Range("D10").Value = "0.0978"
$ws.Range("E10").Value = "  -2.23%  "

# Row 11
$ws.Range("D11").Value = "57.67"
$ws.Range("E11").Value = "  -1.70%  "

# Row 12
$ws.Range("D12").Value = "34.71"
$ws.Range("E12").Value = "  +5.97%  "

# Row 13
$ws.Range("E13").Value = "  -2.96%  "

# Row 14
$ws.Range("D14").Value = "6.69"
$ws.Range("E14").Value = "  -8.39%  "

# Row 15
$ws.Range("D15").Value = "2.567.62"
$ws.Range("E15").Value = "  -4.76%  "

# Row 16
$ws.Range("D16").Value = "14.66"
$ws.Range("E16").Value = "  -10.27%  "

# Row 17
$ws.Range("D17").Value = "0.861"
$ws.Range("E17").Value = "  -4.47%  "

# Row 18
$ws.Range("D18").Value = "2.236.46"
$ws.Range("E18").Value = "  -4.75%  "

# Row 19
$ws.Range("D19").Value = "41.874.64"
$ws.Range("E19").Value = "  -4.16%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -3.85%  "

# Row 21
$ws.Range("D21").Value = "72.95"
$ws.Range("E21").Value = "  -5.59%  "

# Row 22
$ws.Range("D22").Value = "6.15"
$ws.Range("E22").Value = "  -8.50%  "

# Row 23
$ws.Range("D23").Value = "233.55"
$ws.Range("E23").Value = "  -8.88%  "

# Row 24
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("D25").Value = "3.60"
$ws.Range("E25").Value = "  -3.83%  "

# Row 26
$ws.Range("D26").Value = "1.78"
$ws.Range("E26").Value = "  -8.03%  "

# Row 27
$ws.Range("D27").Value = "2.35"
$ws.Range("E27").Value = "  -5.53%  "

# Row 28
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  -6.80%  "

# Row 29
$ws.Range("D29").Value = "2.10"
$ws.Range("E29").Value = "  -7.73%  "

# Row 30
$ws.Range("D30").Value = "165.58"
$ws.Range("E30").Value = "  -6.64%  "

# Row 31
$ws.Range("D31").Value = "20.48"

# Row 32
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  -8.03%  "

# Row 33
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  -8.44%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.0702"
$ws.Range("E34").Value = "  -7.19%  "

# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  -4.66%  "

# Row 36
$ws.Range("D36").Value = "4.75"
$ws.Range("E36").Value = "  -8.18%  "

# Row 37
$ws.Range("D37").Value = "3.55"
$ws.Range("E37").Value = "  -6.89%  "

# Row 38
$ws.Range("D38").Value = "5.98"
$ws.Range("E38").Value = "  -4.75%  "

# Row 39
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  -6.51%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0261"
$ws.Range("E40").Value = "  -6.99%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "20.42"
$ws.Range("E41").Value = "  +6.94%  "

# Row 42
$ws.Range("D42").Value = "63.81"
$ws.Range("E42").Value = "  -8.12%  "

# Row 43
$ws.Range("D43").Value = "5.04"
$ws.Range("E43").Value = "  +4.83%  "

# Row 44
$ws.Range("D44").Value = "8.69"
$ws.Range("E44").Value = "  -4.85%  "

# Row 45
$ws.Range("B45").Value = "BinanceUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0988"
$ws.Range("E46").Value = "  -11.41%  "

# Row 47
$ws.Range("D47").Value = "0.184"
$ws.Range("E47").Value = "  -8.76%  "

# Row 48
$ws.Range("D48").Value = "4.33"
$ws.Range("E48").Value = "  +7.37%  "

# Row 49
$ws.Range("E49").Value = "  -6.39%  "

# Row 50
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  -7.63%  "

# Row 51
$ws.Range("E51").Value = "  -0.98%  "
